# Add team win/loss/tie record columns to the player data sheet.
# New columns: AD = Wins, AE = Losses, AF = Ties.
# Header row (row 1) reuses the existing bold/bordered header style (same
# as A1:AC1); data rows (2-54) get plain numeric values, one row per
# player but all sharing the team's 95-67-0 record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header labels for the three new columns.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the header formatting already used by the rest of row 1 (bold
# font, thin border, centered/top aligned) by copying A1's format instead
# of re-building it by hand, so no new style entries are introduced.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the team record (95 wins, 67 losses, 0 ties) for every player
# row.
$lastRow = 54
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 95
    $ws.Cells.Item($r, 31).Value = 67
    $ws.Cells.Item($r, 32).Value = 0
}
